# Updated symbol list on Sun Dec 25 13:59:57 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking value while keeping it stored as literal
# text (so formatting like trailing zeros, e.g. "244.20", is preserved
# exactly as in the source data) and without leaving a stray custom style
# behind on the cell.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Simple Price (D column) updates that don't involve row moves
Set-TextValue "D2" "244.20"
Set-TextValue "D5" "0.05953"
Set-TextValue "D7" "0.8067"
Set-TextValue "D8" "0.9254"

# Rows 9-17: "One" moved up from row 17 to row 9, and WazirX..CoinExToken
# each shifted down by one row, with several price/volume updates along the way.
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D9" "0.01118"
$ws.Range("E9").Value = "8OneONEBestin24h"

$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1417"
$ws.Range("E10").Value = "9WazirXWRX"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.07424"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"

$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D12" "0.03415"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D13" "0.03070"
$ws.Range("E13").Value = "12BitrueCoinBTR"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D14" "0.09338"
$ws.Range("E14").Value = "13BitMartTokenBMX"

$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D15" "3.939"
$ws.Range("E15").Value = "14MCDexMCB"

$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D16" "0.001595"
$ws.Range("E16").Value = "15BitForexTokenBF"

$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D17" "0.04808"
$ws.Range("E17").Value = "16CoinExTokenCET"

# Remaining simple Price (D column) updates, plus one Volume(1h) (E column) update
Set-TextValue "D18" "0.005640"
Set-TextValue "D19" "0.004154"
Set-TextValue "D20" "0.0009831"
Set-TextValue "D22" "3.654"
Set-TextValue "D23" "6.442"
Set-TextValue "D24" "2.186"
Set-TextValue "D26" "0.1340"
Set-TextValue "D40" "0.03927"
Set-TextValue "D41" "0.006225"
Set-TextValue "D42" "0.1072"
Set-TextValue "D43" "0.002902"
Set-TextValue "D44" "0.007511"
$ws.Range("E44").Value = "43LocalTradersLCT"
Set-TextValue "D45" "0.00005199"
Set-TextValue "D48" "1.050"
Set-TextValue "D51" "0.0002001"
